$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version bump
$ws.Range("B3").Value = "6.0.0"

# Date update
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value
$ws.Range("B9").Value = "Alvearie Team"

# Replace duplicate "Contact" row (row 10) with "Jurisdiction" row
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Remove the now-redundant second "Contact" row (row 11); rows below shift up
$ws.Rows.Item(11).Delete()

# "Case Sensitive" row (now row 14 after the delete) gets a value
$ws.Range("B14").Value = "true"
